$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update length constraints
$ws.Range("B37").Value = "almeno 1 caratteri massimo 30 totali"
$ws.Range("B13").Value = "almeno 3 caratteri massimo 27"
$ws.Range("B15").Value = "almeno 4 caratteri massimo 20"

# Delete rows 38-39 (città / provincia under DD_Ord) - shift rows up
$ws.Range("A38:D39").EntireRow.Delete()
